$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update circuit voltage values (column C) ---
$ws.Range("C2").Value = 127
$ws.Range("C3").Value = 220
$ws.Range("C4").Value = 127

# --- Move the active selection from I16 to J11 ---
$ws.Range("J11").Select()

# --- Best-effort window geometry update (host may not persist these) ---
$win = $excel.ActiveWindow
$win.Width = 21600
$win.Height = 11385
$win.Left = -26985
$win.Top = 1815
